# Generate Report for Handoff
# Update the "Status" columns from "In Translation" to "Ready for handoff"
# and bump the related handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2017-02-22 08:06:05"

# zh-cn sheet: C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2017-02-22 08:05:48"

# de-de sheet: C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2017-02-22 08:06:05"

# Column width adjustments (widened to fit the new "Ready for handoff" text)
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
